$d = $word.ActiveDocument

# --- Define the new character styles ---------------------------------
$gaNStyle = $d.Styles.Add("GaNStyle", 2)
$gaNStyle.Font.Name = "Calibri"
$gaNStyle.Font.Size = 14

$gaNParagraph = $d.Styles.Add("GaNParagraph", 2)
$gaNParagraph.Font.Name = "Calibri"
$gaNParagraph.Font.Size = 10

$gaNLinks = $d.Styles.Add("GaNLinks", 2)
$gaNLinks.Font.Name = "Calibri"
$gaNLinks.Font.Bold = $true
$gaNLinks.Font.Color = 8388608
$gaNLinks.Font.Size = 9.5
$gaNLinks.Font.Underline = 1

# --- Apply GaNStyle to every "Datas das campanhas ..." run -----------
$searchText = "Datas das campanhas de 2022 que usam constelação de leão: 14 a 23 de abril, 14 a 23 de maio"
$rng = $d.Content
$rng.Start = 0
while ($rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNStyle"
    $rng.Collapse(0)
    $rng.End = $d.Content.End
}

# --- Apply GaNParagraph to the "Está a participar ..." run -----------
$searchText2 = "Está a participar numa campanha global para observar e registar as estrelas mais fracas visíveis como forma de medir a poluição luminosa num determinado local. Localizando e observando a  constelação de leão no céu noturno e,  comparando-a com cartas estelares, pessoas de todo o mundo aprenderão  como as luzes da sua comunidade contribuem para a poluição luminosa. As suas contribuições para a base de dados on-line irão documentar a visibilidade do céu noturno em todo o mundo."
$rng2 = $d.Content
$rng2.Start = 0
while ($rng2.Find.Execute($searchText2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng2.Style = "GaNParagraph"
    $rng2.Collapse(0)
    $rng2.End = $d.Content.End
}

# --- Apply GaNLinks to the "por Jenik Hollan ..." run -----------------
$searchText3 = "por Jenik Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."
$rng3 = $d.Content
$rng3.Start = 0
while ($rng3.Find.Execute($searchText3, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng3.Style = "GaNLinks"
    $rng3.Collapse(0)
    $rng3.End = $d.Content.End
}
